# "datos nuevos para testear" - adds new test data (wind/nuclear figures) to the
# "A testear" sheet, removing the never-populated "% that is Private Land" column
# and filling in the previously-empty Difference elevation / % Public Land /
# Installed Wind Capacity / Nuclear Year to Date columns for each country row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A testear")
[void]$ws.Activate()

# The "% that is Private Land" column (H) never had any data in it - drop it so
# the Wind Speed / Wind Capacity / Nuclear headers shift left into H/I/J.
[void]$ws.Columns.Item(8).Delete()

# --- Spain (row 2) ---
$ws.Range("F2").Value = 3718
$ws.Range("G2").Value = 17.1
$ws.Range("I2").Value = 23025
$ws.Range("J2").Value = 54740
$ws.Range("I2:J2").WrapText = $true

# --- France (row 3) ---
$ws.Range("F3").Value = 4810
$ws.Range("G3").Value = 20
$ws.Range("I3").Value = 10358
$ws.Range("J3").Value = 416800
$ws.Range("I3:J3").WrapText = $true

# --- Germany (row 4) ---
$ws.Range("F4").Value = 2962
$ws.Range("G4").Value = 15.4
$ws.Range("I4").Value = 44947
$ws.Range("J4").Value = 86810.32
$ws.Range("J4").WrapText = $true

# --- England (row 5) ---
$ws.Range("F5").Value = 1344
$ws.Range("G5").Value = 23.5
$ws.Range("I5").Value = 13603
$ws.Range("J5").Value = 63894.54
$ws.Range("I5:J5").WrapText = $true

# --- Greece (row 6) ---
$ws.Range("F6").Value = 2919
$ws.Range("G6").Value = 22.6
$ws.Range("I6").Value = 2152
$ws.Range("J6").Value = 0
$ws.Range("I6").WrapText = $true

# The new rows render a touch taller than the sheet's old default.
$ws.Rows.Item(2).RowHeight = 15.65
$ws.Rows.Item(3).RowHeight = 15.65
$ws.Rows.Item(4).RowHeight = 15.65
$ws.Rows.Item(5).RowHeight = 15.65
$ws.Rows.Item(6).RowHeight = 15.65

# Selection ends up parked just past the populated data.
[void]$ws.Range("J7").Select()
